# Updated cryptos list on Fri Apr 28 11:36:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell 2 4 "29.310.64"
Set-TextCell 2 5 "  +0.92%  "
Set-TextCell 3 4 "1.912.59"
Set-TextCell 3 5 "  +1.36%  "
Set-TextCell 4 4 "1.001"
Set-TextCell 4 5 "  -0.05%  "
Set-TextCell 5 4 "321.76"
Set-TextCell 5 5 "  -2.81%  "
Set-TextCell 6 4 "1.001"
Set-TextCell 6 5 "  +0.02%  "
Set-TextCell 7 4 "0.4709"
Set-TextCell 7 5 "  +2.31%  "
Set-TextCell 8 4 "0.4048"
Set-TextCell 8 5 "  -0.49%  "
Set-TextCell 9 2 "OKB"
Set-TextCell 9 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 9 4 "47.68"
Set-TextCell 9 5 "  -0.20%  "
Set-TextCell 10 2 "Dogecoin"
Set-TextCell 10 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 10 4 "0.08032"
Set-TextCell 10 5 "  +0.57%  "
Set-TextCell 11 2 "Polygon"
Set-TextCell 11 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 11 4 "0.9993"
Set-TextCell 11 5 "  +0.67%  "
Set-TextCell 12 2 "Solana"
Set-TextCell 12 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 12 4 "22.63"
Set-TextCell 12 5 "  +4.30%  "
Set-TextCell 13 2 "WrappedEther"
Set-TextCell 13 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 13 4 "1.965.83"
Set-TextCell 13 5 "  +6.96%  "
Set-TextCell 14 2 "Polkadot"
Set-TextCell 14 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 14 4 "5.879"
Set-TextCell 14 5 "  -0.64%  "
Set-TextCell 15 2 "Chainlink"
Set-TextCell 15 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 15 4 "7.104"
Set-TextCell 15 5 "  +0.41%  "
Set-TextCell 16 2 "Litecoin"
Set-TextCell 16 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 16 4 "89.49"
Set-TextCell 16 5 "  +1.07%  "
Set-TextCell 17 2 "BinanceUSD"
Set-TextCell 17 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 17 4 "1.001"
Set-TextCell 17 5 "  -0.10%  "
Set-TextCell 18 2 "TRON"
Set-TextCell 18 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 18 4 "0.06622"
Set-TextCell 18 5 "  +1.03%  "
Set-TextCell 19 2 "ShibaInu"
Set-TextCell 19 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 19 4 "0.00001028"
Set-TextCell 19 5 "  -0.09%  "
Set-TextCell 20 2 "Avalanche"
Set-TextCell 20 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 20 4 "17.63"
Set-TextCell 20 5 "  +1.06%  "
Set-TextCell 21 2 "Dai"
Set-TextCell 21 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 21 4 "1.000"
Set-TextCell 21 5 "  +0.02%  "
Set-TextCell 22 2 "WrappedBTC"
Set-TextCell 22 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 22 4 "29.334.38"
Set-TextCell 22 5 "  +0.91%  "
Set-TextCell 23 2 "Uniswap"
Set-TextCell 23 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 23 4 "5.521"
Set-TextCell 23 5 "  +1.84%  "
Set-TextCell 24 2 "Cosmos"
Set-TextCell 24 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 24 4 "11.40"
Set-TextCell 24 5 "  -0.45%  "
Set-TextCell 25 2 "Toncoin"
Set-TextCell 25 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 25 4 "2.203"
Set-TextCell 25 5 "  -0.47%  "
Set-TextCell 26 2 "WrappedliquidstakedEther2.0"
Set-TextCell 26 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 26 4 "2.139.20"
Set-TextCell 26 5 "  +3.17%  "
Set-TextCell 27 2 "Monero"
Set-TextCell 27 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 27 4 "154.34"
Set-TextCell 27 5 "  -1.71%  "
Set-TextCell 28 2 "EthereumClassic"
Set-TextCell 28 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 28 4 "19.76"
Set-TextCell 28 5 "  +0.85%  "
Set-TextCell 29 2 "InternetComputer(DFINITY)"
Set-TextCell 29 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 29 4 "6.043"
Set-TextCell 29 5 "  +10.44%  "
Set-TextCell 30 2 "LidoDAOToken"
Set-TextCell 30 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 30 4 "2.105"
Set-TextCell 30 5 "  +0.20%  "
Set-TextCell 31 2 "BitcoinCash"
Set-TextCell 31 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 31 4 "117.93"
Set-TextCell 31 5 "  +0.26%  "
Set-TextCell 32 2 "ImmutableX"
Set-TextCell 32 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 32 4 "1.066"
Set-TextCell 32 5 "  +5.49%  "
Set-TextCell 33 2 "Stellar"
Set-TextCell 33 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 33 4 "0.09492"
Set-TextCell 33 5 "  +1.68%  "
Set-TextCell 34 2 "ARBITRUM"
Set-TextCell 34 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 34 4 "1.418"
Set-TextCell 34 5 "  +0.56%  "
Set-TextCell 35 2 "HuobiToken"
Set-TextCell 35 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 35 4 "3.544"
Set-TextCell 35 5 "  -1.57%  "
Set-TextCell 36 2 "Filecoin"
Set-TextCell 36 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 36 4 "5.372"
Set-TextCell 36 5 "  +1.70%  "
Set-TextCell 37 2 "Hedera"
Set-TextCell 37 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 37 4 "0.06074"
Set-TextCell 37 5 "  +0.15%  "
Set-TextCell 38 2 "VeChain"
Set-TextCell 38 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 38 4 "0.02247"
Set-TextCell 38 5 "  +1.06%  "
Set-TextCell 39 2 "FraxShare"
Set-TextCell 39 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 39 4 "8.183"
Set-TextCell 39 5 "  -1.28%  "
Set-TextCell 40 2 "TrustWalletToken"
Set-TextCell 40 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 40 4 "1.177"
Set-TextCell 40 5 "  +0.28%  "
Set-TextCell 41 2 "TheSandbox"
Set-TextCell 41 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 41 4 "0.5845"
Set-TextCell 41 5 "  +0.99%  "
Set-TextCell 42 2 "Algorand"
Set-TextCell 42 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 42 4 "0.1836"
Set-TextCell 42 5 "  +0.71%  "
Set-TextCell 43 2 "RenderToken"
Set-TextCell 43 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 43 4 "2.492"
Set-TextCell 43 5 "  +9.70%  "
Set-TextCell 44 2 "Aptos"
Set-TextCell 44 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 44 4 "10.10"
Set-TextCell 44 5 "  -0.41%  "
Set-TextCell 45 2 "Cronos"
Set-TextCell 45 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell 45 4 "0.07911"
Set-TextCell 45 5 "  +5.67%  "
Set-TextCell 46 2 "WEMIXToken"
Set-TextCell 46 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 46 4 "1.277"
Set-TextCell 46 5 "  +1.43%  "
Set-TextCell 47 4 "0.5503"
Set-TextCell 47 5 "  +1.04%  "
Set-TextCell 48 2 "EnergySwap"
Set-TextCell 48 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 48 4 "12.05"
Set-TextCell 48 5 "  -0.03%  "
Set-TextCell 49 2 "NEARProtocol"
Set-TextCell 49 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 49 4 "1.923"
Set-TextCell 49 5 "  +1.11%  "
Set-TextCell 50 2 "Quant"
Set-TextCell 50 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell 50 4 "113.11"
Set-TextCell 50 5 "  +1.82%  "
Set-TextCell 51 2 "Elrond"
Set-TextCell 51 3 "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell 51 4 "44.21"
Set-TextCell 51 5 "  -2.97%  "
